$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 (Japanese column) previously held a placeholder duplicate of the
# Chinese text; replace it with the real Japanese translation.
$ws.Range("B2").Value = "アンジェリーナは、荷物の受取人を見つけられなかった。しかしロドスに戻ると、その荷物の差出人がスズランで、受取人が自分だったことをようやく知るのだった。
"

# C2 (English column) previously held the same placeholder duplicate of
# the Chinese text; replace it with the real English translation.
$ws.Range("C2").Value = "Unable to find the package's owner, Angelina returns to Rhodes Island only to find that the owner of the package in her hand is Suzuran.
"
